# Generate Report for Handoff
#
# - Status moves from "In Translation" to "Ready for handoff" everywhere it
#   appears: Overview!E2:F2, zh-cn!C2, de-de!C2.
# - The handoff timestamps tied to that status move forward too:
#     Overview!G2 (Latest HO Xliff Generate Date)  10:54:52 -> 10:55:28
#     zh-cn!H2    (Latest Handoff Datetime)         10:54:47 -> 10:55:23
#     de-de!H2    (Latest Handoff Datetime)         10:54:52 -> 10:55:28
# - The Status columns (now holding the longer "Ready for handoff" text)
#   widen to fit.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item(1)
$zhcn     = $wb.Worksheets.Item(2)
$dede     = $wb.Worksheets.Item(3)

# --- Status text: "In Translation" -> "Ready for handoff" ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value     = "Ready for handoff"
$dede.Range("C2").Value     = "Ready for handoff"

# --- Handoff timestamps bumped alongside the status change ---
$overview.Range("G2").Value = "2016-08-17 10:55:28"
$zhcn.Range("H2").Value     = "2016-08-17 10:55:23"
$dede.Range("H2").Value     = "2016-08-17 10:55:28"

# --- Widen the Status columns so the longer text still fits ---
$overview.Columns.Item(5).ColumnWidth = 16.333333333333336
$overview.Columns.Item(6).ColumnWidth = 16.333333333333336
$zhcn.Columns.Item(3).ColumnWidth     = 16.333333333333336
$dede.Columns.Item(3).ColumnWidth     = 16.333333333333336
